$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-10 from serial 45174 to 45175
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}
